$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.162165
$ws.Range("N2").Value = 0.486495
$ws.Range("O2").Value = 0.1006291402646046
$ws.Range("P2").Value = 0.1006291402646046
$ws.Range("Q2").Value = 0.01155706711
$ws.Range("R2").Value = 0.10401360399
$ws.Range("S2").Value = 0.1006291402646046
$ws.Range("T2").Value = 0.1006291402646046

# Row 3
$ws.Range("O3").Value = 0.1908661724170313
$ws.Range("P3").Value = 0.1908661724170313
$ws.Range("S3").Value = 0.1908661724170313
$ws.Range("T3").Value = 0.1908661724170313

# Row 4
$ws.Range("M4").Value = 1.141763333333333
$ws.Range("N4").Value = 3.42529
$ws.Range("O4").Value = 0.7085046873183641
$ws.Range("P4").Value = 0.7085046873183641
$ws.Range("Q4").Value = 0.08137042806444446
$ws.Range("R4").Value = 0.7323338525800001
$ws.Range("S4").Value = 0.7085046873183641
$ws.Range("T4").Value = 0.7085046873183641
